# Fix champion names and convert ASSISTS column (F) from text to numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F (ASSISTS) rows 2..41 were stored as text ("0","1","2") - convert to
# real numbers while keeping the same displayed value.
for ($r = 2; $r -le 41; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $cell.Value = [double]$cell.Text
}

# Column H (CHAMPION): several rows were mis-labeled and should all read
# "Smolder".
$rowsToFix = @(5, 11, 17, 23, 29, 30, 35, 36, 41)
foreach ($r in $rowsToFix) {
    $ws.Cells.Item($r, 8).Value = "Smolder"
}
